# Machine_Service_Lookup.xlsx - Card24 sheet
# Add a new service-range row under the "0-150" range (admin manually added
# a calibration record), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Insert a new blank row at position 3 (old rows 3..12 shift to 4..13)
$ws.Rows.Item(3).Insert()

# The numeric-looking columns (card/Min_Tones/Max_Tones/Tones) in this sheet
# are stored as text everywhere else in the table, so format the new cells
# as text before writing them to avoid Excel auto-converting them to numbers.
$ws.Range("A3:D3").NumberFormat = "@"

$ws.Range("A3").Value = "24"
$ws.Range("B3").Value = "0"
$ws.Range("C3").Value = "150"
$ws.Range("D3").Value = "99"
$ws.Range("E3").Value = "done"
$ws.Range("F3").Value = "done"
$ws.Range("M3").Value = " تم معايره "

# G3,H3,I3,J3,K3,L3,N3 are left blank for this new row.
